$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("14:17").Delete()

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Mdk"
$ws.Range("C2").Value = "Tspan1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.988074333333333
$ws.Range("H2").Value = 5.964223
$ws.Range("I2").Value = 0.01657769708907969
$ws.Range("J2").Value = 0.01657769708907968
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.583904
$ws.Range("N2").Value = 1.751712
$ws.Range("O2").Value = 0.2526535508491896
$ws.Range("P2").Value = 0.2526535508491896
$ws.Range("Q2").Value = 1.160844555530667
$ws.Range("R2").Value = 10.447600999776
$ws.Range("S2").Value = 0.004188414034458257
$ws.Range("T2").Value = 0.004188414034458256

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Mdk"
$ws.Range("C3").Value = "Tspan1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.988074333333333
$ws.Range("H3").Value = 5.964223
$ws.Range("I3").Value = 0.01657769708907969
$ws.Range("J3").Value = 0.01657769708907968
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.9102196666666668
$ws.Range("N3").Value = 2.730659
$ws.Range("O3").Value = 0.3938493842071626
$ws.Range("P3").Value = 0.3938493842071626
$ws.Range("Q3").Value = 1.809584356995222
$ws.Range("R3").Value = 16.286259212957
$ws.Range("S3").Value = 0.006529115790106907
$ws.Range("T3").Value = 0.006529115790106906

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Mdk"
$ws.Range("C4").Value = "Tspan1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.988074333333333
$ws.Range("H4").Value = 5.964223
$ws.Range("I4").Value = 0.01657769708907969
$ws.Range("J4").Value = 0.01657769708907968
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.8169620000000001
$ws.Range("N4").Value = 2.450886
$ws.Range("O4").Value = 0.3534970649436477
$ws.Range("P4").Value = 0.3534970649436476
$ws.Range("Q4").Value = 1.624181183508667
$ws.Range("R4").Value = 14.617630651578
$ws.Range("S4").Value = 0.005860167264514522
$ws.Range("T4").Value = 0.005860167264514519

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Mdk"
$ws.Range("C5").Value = "Tspan1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 92.89399466666667
$ws.Range("H5").Value = 278.681984
$ws.Range("I5").Value = 0.7746030815641455
$ws.Range("J5").Value = 0.7746030815641454
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.583904
$ws.Range("N5").Value = 1.751712
$ws.Range("O5").Value = 0.2526535508491896
$ws.Range("P5").Value = 0.2526535508491896
$ws.Range("Q5").Value = 54.24117506184533
$ws.Range("R5").Value = 488.170575556608
$ws.Range("S5").Value = 0.1957062190559058
$ws.Range("T5").Value = 0.1957062190559058

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Mdk"
$ws.Range("C6").Value = "Tspan1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 92.89399466666667
$ws.Range("H6").Value = 278.681984
$ws.Range("I6").Value = 0.7746030815641455
$ws.Range("J6").Value = 0.7746030815641454
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.9102196666666668
$ws.Range("N6").Value = 2.730659
$ws.Range("O6").Value = 0.3938493842071626
$ws.Range("P6").Value = 0.3938493842071626
$ws.Range("Q6").Value = 84.55394086082846
$ws.Range("R6").Value = 760.985467747456
$ws.Range("S6").Value = 0.3050769466790093
$ws.Range("T6").Value = 0.3050769466790093

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Mdk"
$ws.Range("C7").Value = "Tspan1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 92.89399466666667
$ws.Range("H7").Value = 278.681984
$ws.Range("I7").Value = 0.7746030815641455
$ws.Range("J7").Value = 0.7746030815641454
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.8169620000000001
$ws.Range("N7").Value = 2.450886
$ws.Range("O7").Value = 0.3534970649436477
$ws.Range("P7").Value = 0.3534970649436476
$ws.Range("Q7").Value = 75.89086367086935
$ws.Range("R7").Value = 683.017773037824
$ws.Range("S7").Value = 0.2738199158292304
$ws.Range("T7").Value = 0.2738199158292303

$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Mdk"
$ws.Range("C8").Value = "Tspan1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 23.741365
$ws.Range("H8").Value = 71.224095
$ws.Range("I8").Value = 0.1979690350870239
$ws.Range("J8").Value = 0.1979690350870239
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.583904
$ws.Range("N8").Value = 1.751712
$ws.Range("O8").Value = 0.2526535508491896
$ws.Range("P8").Value = 0.2526535508491896
$ws.Range("Q8").Value = 13.86267798896
$ws.Range("R8").Value = 124.76410190064
$ws.Range("S8").Value = 0.0500175796729244
$ws.Range("T8").Value = 0.0500175796729244

$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Mdk"
$ws.Range("C9").Value = "Tspan1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 23.741365
$ws.Range("H9").Value = 71.224095
$ws.Range("I9").Value = 0.1979690350870239
$ws.Range("J9").Value = 0.1979690350870239
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.9102196666666668
$ws.Range("N9").Value = 2.730659
$ws.Range("O9").Value = 0.3938493842071626
$ws.Range("P9").Value = 0.3938493842071626
$ws.Range("Q9").Value = 21.60985733651167
$ws.Range("R9").Value = 194.488716028605
$ws.Range("S9").Value = 0.07796998256111055
$ws.Range("T9").Value = 0.07796998256111055

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Mdk"
$ws.Range("C10").Value = "Tspan1"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 23.741365
$ws.Range("H10").Value = 71.224095
$ws.Range("I10").Value = 0.1979690350870239
$ws.Range("J10").Value = 0.1979690350870239
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.8169620000000001
$ws.Range("N10").Value = 2.450886
$ws.Range("O10").Value = 0.3534970649436477
$ws.Range("P10").Value = 0.3534970649436476
$ws.Range("Q10").Value = 19.39579303313
$ws.Range("R10").Value = 174.56213729817
$ws.Range("S10").Value = 0.06998147285298896
$ws.Range("T10").Value = 0.06998147285298896

$ws.Range("A11").Value = "Resolving-Mac"
$ws.Range("B11").Value = "Mdk"
$ws.Range("C11").Value = "Tspan1"
$ws.Range("D11").Value = "ECs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.301204666666667
$ws.Range("H11").Value = 3.903614
$ws.Range("I11").Value = 0.01085018625975097
$ws.Range("J11").Value = 0.01085018625975097
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.583904
$ws.Range("N11").Value = 1.751712
$ws.Range("O11").Value = 0.2526535508491896
$ws.Range("P11").Value = 0.2526535508491896
$ws.Range("Q11").Value = 0.7597786096853334
$ws.Range("R11").Value = 6.838007487168
$ws.Range("S11").Value = 0.002741338085901171
$ws.Range("T11").Value = 0.00274133808590117

$ws.Range("A12").Value = "Resolving-Mac"
$ws.Range("B12").Value = "Mdk"
$ws.Range("C12").Value = "Tspan1"
$ws.Range("D12").Value = "FAPs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.301204666666667
$ws.Range("H12").Value = 3.903614
$ws.Range("I12").Value = 0.01085018625975097
$ws.Range("J12").Value = 0.01085018625975097
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.9102196666666668
$ws.Range("N12").Value = 2.730659
$ws.Range("O12").Value = 0.3938493842071626
$ws.Range("P12").Value = 0.3938493842071626
$ws.Range("Q12").Value = 1.184382077958445
$ws.Range("R12").Value = 10.659438701626
$ws.Range("S12").Value = 0.004273339176935937
$ws.Range("T12").Value = 0.004273339176935937

$ws.Range("A13").Value = "Resolving-Mac"
$ws.Range("B13").Value = "Mdk"
$ws.Range("C13").Value = "Tspan1"
$ws.Range("D13").Value = "MuSCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.301204666666667
$ws.Range("H13").Value = 3.903614
$ws.Range("I13").Value = 0.01085018625975097
$ws.Range("J13").Value = 0.01085018625975097
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.8169620000000001
$ws.Range("N13").Value = 2.450886
$ws.Range("O13").Value = 0.3534970649436477
$ws.Range("P13").Value = 0.3534970649436476
$ws.Range("Q13").Value = 1.063034766889334
$ws.Range("R13").Value = 9.567312902004002
$ws.Range("S13").Value = 0.003835508996913863
$ws.Range("T13").Value = 0.003835508996913862
